# Auto-generated edit script: updates crypto price/volume table
# to match the target diff (commit: "Updated cryptos list ... with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "28.795.99"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  +7.26%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.813.36"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  +5.07%  "

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "0.9994"
$cell.Style = "Normal"
$ws.Range("E4").Value = "  +0.19%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "250.65"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +3.61%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "0.9994"
$cell.Style = "Normal"

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.4985"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  +2.05%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.2774"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  +7.15%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.06379"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +2.75%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "1.813.63"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +5.04%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "16.76"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  +4.90%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.07156"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  +3.70%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "0.6488"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +6.90%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "4.708"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  +5.23%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "81.96"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +6.17%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "28.776.57"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +8.07%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "0.9987"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +0.10%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "0.000007377"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +2.82%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "0.9994"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +0.17%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "12.27"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +7.13%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "2.047.12"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +4.89%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "4.610"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +4.45%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "8.887"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +3.64%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "5.346"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +5.20%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "143.78"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +3.87%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "16.01"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +4.61%  "

$ws.Range("B27").Value = "BitcoinCash"
$ws.Range("C27").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "115.87"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  +8.99%  "

$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "1.888"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +6.59%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "1.394"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +1.05%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "4.174"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +5.76%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "0.08358"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  +4.73%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "3.844"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  +4.40%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "0.04973"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +10.04%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "1.089"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +7.98%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "0.6781"
$cell.Style = "Normal"

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "2.683"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +3.25%  "

$ws.Range("E37").Value = "  +11.97%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.9697"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +3.85%  "

$ws.Range("E39").Value = "  +7.05%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.01591"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +6.25%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "5.997"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +5.74%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.9993"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  +0.18%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "101.64"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +2.16%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.4116"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +7.20%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "7.227"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  +5.55%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.1226"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  +5.62%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "0.05505"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +1.96%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "8.181"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +3.53%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "31.65"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +4.85%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "0.3651"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  +8.32%  "

$ws.Range("E51").Value = "  +5.88%  "

